$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(6, 0.2590361445783131, 0.2453109174966812, 0.01372522708163193, 0.9864606042002675)
    3 = @(5, 0.2048192771084338, 0.2393135130405426, 0.03449423593210885, 0.9666559418757626)
    4 = @(0, 0.1927710843373494, 0.2334723174571991, 0.04070123311984972, 0.9608905689504814)
    5 = @(3, 0.9156626506024095, 0.9610681533813477, 0.04540550277893818, 0.9565666120388313)
    6 = @(7, 0.06265060240963854, 0.1267924308776855, 0.064141828468047, 0.9397243612156601)
    7 = @(2, 0.1349397590361446, 0.2219508737325668, 0.08701111469642228, 0.9199537948416263)
    8 = @(1, 0.03855421686746985, 0.1274387687444687, 0.08888455187699884, 0.918371004783031)
    9 = @(4, 0.05301204819277108, 0.1701029539108276, 0.1170909057180566, 0.8951822943695066)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
}
